{"js": "// Apply the \"Added many more features\" edits to the E-Force review doc.\n// Each entry is an exact, case-sensitive search string mapped to its\n// replacement. Using body.search keeps this robust to paragraph/run\n// splitting quirks, and covers the title text which appears twice\n// (the Heading1 at the top and the bold recap line near the end).\nconst replacements = [\n  [\n    \"Play E-Force Slot Game for Free - Review\",\n    \"Play E-Force Slot Game for Free\"\n  ],\n  [\n    \"Humorous animations and cartoonish graphics\",\n    \"Visually appealing graphics with a cartoonish style\"\n  ],\n  [\n    \"Catchy tune played on real instruments\",\n    \"Humorous animations enhance the gameplay\"\n  ],\n  [\n    \"243 ways to win with a betting range of \\u20AC0.10 to \\u20AC100\",\n    \"Catchy tune played with real instruments adds to the overall experience\"\n  ],\n  [\n    \"Players can purchase the free spin feature\",\n    \"Exciting free spin feature with the opportunity to earn additional spins and increasing multipliers\"\n  ],\n  [\n    \"The theme may not be appealing to all players\",\n    \"Limited number of high-level symbols\"\n  ],\n  [\n    \"No progressive jackpot\",\n    \"Purchasing the free spin feature can be expensive\"\n  ],\n  [\n    \"Read our E-Force slot game review and play for free. Learn about the gameplay, symbols, and free spin feature in this Yggdrasil-developed game.\",\n    \"Read our review of E-Force, a slot game developed by Yggdrasil and play for free to experience the exciting free spin feature.\"\n  ]\n];\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the E-Force review doc.\n# Uses Find/Replace (wdReplaceAll) over the whole document Range for each\n# exact phrase. The title phrase appears twice (Heading1 at the top, and\n# the bold recap line near the end) - ReplaceAll handles both in one call.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute(\n        $findText,    # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n}\n\nReplace-Text \"Play E-Force Slot Game for Free - Review\" \"Play E-Force Slot Game for Free\"\n\nReplace-Text \"Humorous animations and cartoonish graphics\" \"Visually appealing graphics with a cartoonish style\"\nReplace-Text \"Catchy tune played on real instruments\" \"Humorous animations enhance the gameplay\"\nReplace-Text \"243 ways to win with a betting range of \u20ac0.10 to \u20ac100\" \"Catchy tune played with real instruments adds to the overall experience\"\nReplace-Text \"Players can purchase the free spin feature\" \"Exciting free spin feature with the opportunity to earn additional spins and increasing multipliers\"\n\nReplace-Text \"The theme may not be appealing to all players\" \"Limited number of high-level symbols\"\nReplace-Text \"No progressive jackpot\" \"Purchasing the free spin feature can be expensive\"\n\nReplace-Text \"Read our E-Force slot game review and play for free. Learn about the gameplay, symbols, and free spin feature in this Yggdrasil-developed game.\" \"Read our review of E-Force, a slot game developed by Yggdrasil and play for free to experience the exciting free spin feature.\"\n"}
